# Update workbook with correct forecast output:
#  - rename Sheet1 -> "Sales vs PO"
#  - add "Weekly Growth", "Volume Insights", "Prediction Info" sheets
#  - insert an "Order Week" column (old ds values) before PO_Requested_Qty
#    on the Sales vs PO sheet, shift ds forward a week, and zero out the
#    PO_Requested_Qty column there
#  - move the non-zero historical PO qty rows (with week-over-week growth%)
#    onto the new Weekly Growth sheet
#  - summarise totals on Volume Insights
#  - seed the (zero, for now) forecast on Prediction Info

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheets
# ---------------------------------------------------------------------
$wsSales = $wb.Worksheets.Item(1)
$wsSales.Name = "Sales vs PO"

$wsGrowth = $wb.Worksheets.Add($null, $wsSales)
$wsGrowth.Name = "Weekly Growth"

$wsVolume = $wb.Worksheets.Add($null, $wsGrowth)
$wsVolume.Name = "Volume Insights"

$wsPred = $wb.Worksheets.Add($null, $wsVolume)
$wsPred.Name = "Prediction Info"

# ---------------------------------------------------------------------
# Sales vs PO: ds | y | Order Week | PO_Requested_Qty
# ---------------------------------------------------------------------
$wsSales.Range("C1").Value = "Order Week"
$wsSales.Range("D1").Value = "PO_Requested_Qty"

# copy the bold/bordered header style onto the two new header cells
$wsSales.Range("A1").Copy()
$wsSales.Range("C1:D1").PasteSpecial(-4122)

$salesData = @(
  @(45557,0,45551,0),
  @(45564,0,45558,0),
  @(45571,0,45565,0),
  @(45578,15,45572,0),
  @(45585,18,45579,0),
  @(45592,17,45586,0),
  @(45599,54,45593,0),
  @(45606,54,45600,0),
  @(45613,89,45607,0),
  @(45620,46,45614,0),
  @(45627,23,45621,0),
  @(45634,108,45628,0),
  @(45641,68,45635,0),
  @(45648,51,45642,0),
  @(45655,73,45649,0)
)

$r = 2
foreach ($row in $salesData) {
    $wsSales.Cells.Item($r, 1).Value = $row[0]
    $wsSales.Cells.Item($r, 2).Value = $row[1]
    $wsSales.Cells.Item($r, 3).Value = $row[2]
    $wsSales.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}

# copy the date style from the old ds column onto the new ds + Order Week columns
$wsSales.Range("A2").Copy()
$wsSales.Range("A2:A16").PasteSpecial(-4122)
$wsSales.Range("C2:C16").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# Weekly Growth: ds | PO_Requested_Qty | Growth%
# ---------------------------------------------------------------------
$wsGrowth.Range("A1").Value = "ds"
$wsGrowth.Range("B1").Value = "PO_Requested_Qty"
$wsGrowth.Range("C1").Value = "Growth%"

$wsSales.Range("A1").Copy()
$wsGrowth.Range("A1:C1").PasteSpecial(-4122)

$growthData = @(
  @(45558,580,0),
  @(45572,20,-96.55172413793103),
  @(45586,20,0),
  @(45593,100,400),
  @(45600,30,-70)
)

$r = 2
foreach ($row in $growthData) {
    $wsGrowth.Cells.Item($r, 1).Value = $row[0]
    $wsGrowth.Cells.Item($r, 2).Value = $row[1]
    $wsGrowth.Cells.Item($r, 3).Value = $row[2]
    $r = $r + 1
}

$wsSales.Range("A2").Copy()
$wsGrowth.Range("A2:A6").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# Volume Insights: Total_PO_Quantity | Average_PO_Quantity | Max_PO_Quantity | Min_PO_Quantity
# ---------------------------------------------------------------------
$wsVolume.Range("A1").Value = "Total_PO_Quantity"
$wsVolume.Range("B1").Value = "Average_PO_Quantity"
$wsVolume.Range("C1").Value = "Max_PO_Quantity"
$wsVolume.Range("D1").Value = "Min_PO_Quantity"

$wsSales.Range("A1").Copy()
$wsVolume.Range("A1:D1").PasteSpecial(-4122)

$wsVolume.Range("A2").Value = 750
$wsVolume.Range("B2").Value = 150
$wsVolume.Range("C2").Value = 580
$wsVolume.Range("D2").Value = 20

# ---------------------------------------------------------------------
# Prediction Info: Predicted_Next_Week_PO_Quantity
# ---------------------------------------------------------------------
$wsPred.Range("A1").Value = "Predicted_Next_Week_PO_Quantity"

$wsSales.Range("A1").Copy()
$wsPred.Range("A1").PasteSpecial(-4122)

$wsPred.Range("A2").Value = 0

# ---------------------------------------------------------------------
# leave the workbook selection back on the first sheet, first cell
# ---------------------------------------------------------------------
$wsSales.Range("A1").Select()
$wsSales.Activate()
